$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LJ Speech")

# Row 2: B2 "<all>" -> "<or>", C2 19 -> 18
$ws.Range("B2").Value = "<or>"
$ws.Range("C2").Value = 18

# Row 3: C3 20 -> 16
$ws.Range("C3").Value = 16

# Row 4: C4 19 -> 11
$ws.Range("C4").Value = 11

# Row 5: C5 20 -> 15
$ws.Range("C5").Value = 15

# Row 6: C6 15 -> 14
$ws.Range("C6").Value = 14

# Row 7: C7 10 -> 9
$ws.Range("C7").Value = 9

# Row 8: C8 19 -> 18
$ws.Range("C8").Value = 18

# Row 9: C9 13 -> 14
$ws.Range("C9").Value = 14

# Row 10: C10 13 -> 11
$ws.Range("C10").Value = 11

# Row 11: C11 20 -> 15
$ws.Range("C11").Value = 15

# Row 12: B12 "<him>" -> "<min>", C12 12 -> 13
$ws.Range("B12").Value = "<min>"
$ws.Range("C12").Value = 13

# Row 13: C13 19 -> 15
$ws.Range("C13").Value = 15

# Row 14: B14 "<all>" -> "<are>", C14 12 -> 7
$ws.Range("B14").Value = "<are>"
$ws.Range("C14").Value = 7

# Row 15: C15 9 -> 10
$ws.Range("C15").Value = 10

# Row 16: C16 10 -> 8
$ws.Range("C16").Value = 8

# Row 17: C17 17 -> 13
$ws.Range("C17").Value = 13

# Row 18: C18 13 -> 12
$ws.Range("C18").Value = 12
